# Apply the record re-shuffle described by the diff: several rows in the
# "Artfynd" sheet had their data swapped/rotated with neighboring rows
# (the row *position* in the sheet stayed the same, but the record that
# lives there moved to another row). We replicate this by swapping the
# underlying cell values between the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- helpers ---------------------------------------------------------

function Get-CellRaw {
    param($ws, $row, $col)
    return $ws.Cells.Item($row, $col).Value2
}

function Set-CellRaw {
    param($ws, $row, $col, $val)
    $cell = $ws.Cells.Item($row, $col)
    if ($null -eq $val) {
        $cell.Value2 = $null
        return
    }
    # Guard against Excel's automatic "this looks like a date" conversion:
    # several text columns store plain ISO date strings (e.g. "2026-01-25")
    # as literal text, not real dates. Assigning such a string straight
    # through Value2 makes Excel re-interpret it as a date serial. Prefixing
    # with an apostrophe (the normal Excel "force text" input) keeps it text.
    if ($val -is [string] -and $val -match '^\d{4}-\d{2}-\d{2}$') {
        $cell.Value2 = "'" + $val
    } else {
        $cell.Value2 = $val
    }
}

# Column span used by the sheet (A .. AY == 1 .. 51)
$FirstCol = 1
$LastCol = 51

# Swap just a specific set of columns (by 1-based column index) between two rows.
function Swap-Cols {
    param($ws, $row1, $row2, [int[]]$cols)
    foreach ($c in $cols) {
        $v1 = Get-CellRaw $ws $row1 $c
        $v2 = Get-CellRaw $ws $row2 $c
        Set-CellRaw $ws $row1 $c $v2
        Set-CellRaw $ws $row2 $c $v1
    }
}

# Rotate a specific set of columns across an ordered list of rows:
# new(row[i]) = old(row[i+1]), wrapping around (matches the diff's 3-way cycle).
function Rotate-Cols {
    param($ws, [int[]]$rows, [int[]]$cols)
    $n = $rows.Count
    foreach ($c in $cols) {
        $olds = New-Object System.Collections.ArrayList
        foreach ($r in $rows) {
            [void]$olds.Add((Get-CellRaw $ws $r $c))
        }
        for ($i = 0; $i -lt $n; $i++) {
            $src = $olds[($i + 1) % $n]
            Set-CellRaw $ws $rows[$i] $c $src
        }
    }
}

# Swap an entire row's contents (all columns) with another row, keeping each
# row's own row-number (so only the data moves, not the row itself).
function Swap-FullRows {
    param($ws, $row1, $row2)
    Swap-Cols $ws $row1 $row2 @($FirstCol..$LastCol)
}

# ---- apply the edits ---------------------------------------------------

# Rows 9 <-> 10: Id/Ost/Nord (A, Q, R) swapped.
Swap-Cols $ws 9 10 @(1, 17, 18)

# Rows 15 -> 16 -> 17 -> 15: Id/Ost/Nord (A, Q, R) rotated.
Rotate-Cols $ws @(15, 16, 17) @(1, 17, 18)

# Rows 25 <-> 26: entire records swapped.
Swap-FullRows $ws 25 26

# Rows 32 <-> 33: entire records swapped.
Swap-FullRows $ws 32 33

# Rows 38 <-> 39: entire records swapped.
Swap-FullRows $ws 38 39

# Rows 52 <-> 53: Id/Ost/Nord/Starttid/Sluttid (A, Q, R, Z, AB) swapped.
Swap-Cols $ws 52 53 @(1, 17, 18, 26, 28)

# Rows 56 <-> 58: Id/Ost/Nord (A, Q, R) swapped.
Swap-Cols $ws 56 58 @(1, 17, 18)

# Rows 57 <-> 59: entire records swapped.
Swap-FullRows $ws 57 59
